$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: increase height from 50 to 60 ---
$ws.Rows.Item(1).RowHeight = 60

# --- Populate new row 5 by copying formatting+values from row 2, then overwrite ---
$ws.Range("A2:T2").Copy($ws.Range("A5:T5"))
$ws.Rows.Item(5).RowHeight = 180

# --- Row 2: update values (id 111 -> 138, English/test -> Math/看图填数 question) ---
$ws.Cells.Item(2,1).Value = 138
$ws.Cells.Item(2,3).Value = "数学"
$ws.Cells.Item(2,4).Value = "看图填数"
$ws.Cells.Item(2,5).Value = "(1) 小动物们排成一排，从右边数，小猴子排在第3个，这一排小动物共有 (        ) 个，被挡住的小动物有 (        ) 个。`n(2) 最左边的2个小动物排到队友的最右边后，从左边数，小猴子排在第 (        ) 个，它回家后，还有 (        ) 个小动物。"
$ws.Cells.Item(2,7).Value = "图片下方描述111"
$ws.Cells.Item(2,11).Value = "看图填数"
$ws.Cells.Item(2,13).Value = "小学数学单元过关练习"
$ws.Cells.Item(2,14).Value = 124
$ws.Cells.Item(2,15).Value = "daf2fe1cc1b294de547d7e45ad932567"
$ws.Cells.Item(2,16).Value = "4,5,6"
$ws.Cells.Item(2,17).Value = $true

# --- Row 3: update id, G cleared, N/O/P/S updated ---
$ws.Cells.Item(3,1).Value = 139
$ws.Cells.Item(3,7).ClearContents()
$ws.Cells.Item(3,14).Value = 51
$ws.Cells.Item(3,15).Value = "7c7972d70c830a9300b2f3be0c838a72"
$ws.Cells.Item(3,16).Value = "1,2,3"
$ws.Cells.Item(3,19).Value = 45223

# --- Row 4: id, title/body replaced with test123, K/M cleared, N/O/S updated ---
$ws.Cells.Item(4,1).Value = 140
$ws.Cells.Item(4,4).Value = "test123"
$ws.Cells.Item(4,5).ClearContents()
$ws.Cells.Item(4,11).Value = "未分类"
$ws.Cells.Item(4,13).ClearContents()
$ws.Cells.Item(4,14).Value = 24
$ws.Cells.Item(4,15).Value = "d41d8cd98f00b204e9800998ecf8427e"
$ws.Cells.Item(4,16).ClearContents()
$ws.Cells.Item(4,19).Value = 45223

# --- Row 5: brand-new question row (English/test, id 141) ---
$ws.Cells.Item(5,1).Value = 141
$ws.Cells.Item(5,2).Value = 1
$ws.Cells.Item(5,3).Value = "英语"
$ws.Cells.Item(5,4).Value = "test"
$ws.Cells.Item(5,5).Value = "123"
$ws.Cells.Item(5,7).ClearContents()
$ws.Cells.Item(5,8).Value = $true
$ws.Cells.Item(5,9).Value = 1
$ws.Cells.Item(5,10).Value = 1
$ws.Cells.Item(5,11).Value = "未分类"
$ws.Cells.Item(5,12).Value = $true
$ws.Cells.Item(5,13).ClearContents()
$ws.Cells.Item(5,14).Value = 0
$ws.Cells.Item(5,15).Value = "8f2acdcbe52f06a042904069ef9310a8"
$ws.Cells.Item(5,16).Value = "1,2,3"
$ws.Cells.Item(5,17).Value = $false
$ws.Cells.Item(5,18).Value = "ben"
$ws.Cells.Item(5,19).Value = 45223
$ws.Cells.Item(5,20).Value = 45223

# --- Move the second picture (previously anchored at row 4/1-indexed = row index 3)
# down to the newly inserted row 5 (row index 4, 0-indexed), following its
# originally-described content which now lives in row 3 of the new layout.
$picture = $ws.Shapes.Item(2)
$picture.Top = $ws.Rows.Item(5).Top
